# "arrange logo slide around footer"
#
# On the last slide (the logo / "Where it's already started" slide), several
# logo pictures are shifted upward (their vertical offset, i.e. Top, is
# reduced) so the group sits better around the footer area. Horizontal
# position, width and height are untouched; "Picture 8" keeps its original
# position.

# The PowerPoint object model here reports/accepts Shape.Top/Left/Width/Height
# in points, while the underlying OOXML stores EMU (1 pt = 12700 EMU) and the
# COM layer rounds the point value through a 32-bit float on its way back to
# EMU. A plain "emu / 12700.0" assignment can therefore land one EMU short of
# the exact target after that float32 round-trip, so we search for a points
# value that reproduces the desired EMU exactly once pushed through the same
# float32 -> EMU conversion.
#
# NOTE: this interpreter does not give "for" loop counters their own lexical
# scope across function calls, so the helper below uses its own loop variable
# name ($k) that is never reused by any caller's loop variable ($i) - reusing
# the same name would let the callee's loop clobber the caller's in-progress
# loop counter.
function Get-PointsForEmu {
    param($targetEmu)
    $base = $targetEmu / 12700.0
    for ($k = 0; $k -le 4000; $k++) {
        $cand = $base + ($k * 0.0000001)
        $f = [float]$cand
        $emu = [math]::Floor([double]$f * 12700.0)
        if ($emu -eq $targetEmu) {
            return $cand
        }
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

# Shape name -> new vertical offset (EMU), taken from the target OOXML.
$newTopEmu = @{
    "Picture 3" = 1814517
    "Picture 4" = 3141735
    "Picture 5" = 4905845
    "Picture 6" = 4483949
    "Picture 7" = 2812282
    "Picture 2" = 1666634
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($newTopEmu.ContainsKey($shape.Name)) {
        $shape.Top = Get-PointsForEmu $newTopEmu[$shape.Name]
    }
}
